$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting the existing rows 9-125 down to 10-126
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the new row 9 (a new weekly sample). D (Fecha) and J (Volumen) carry the
# new observation; the remaining fields mirror the template row as in the source data.
$ws.Range('A9').Value = 10
$ws.Range('B9').Value = 'Vega Modelo de Temuco'
$ws.Range('C9').Value = 'La Araucanía'
$ws.Range('D9').Value = 45169
$ws.Range('E9').Value = 9
$ws.Range('F9').Value = 300000001
$ws.Range('G9').Value = 'Rabanito'
$ws.Range('H9').Value = 'Sin especificar'
$ws.Range('I9').Value = 'Primera'
$ws.Range('J9').Value = 80
$ws.Range('K9').Value = 8000
$ws.Range('L9').Value = 8000
$ws.Range('M9').Value = 8000
$ws.Range('N9').Value = '$/docena de paquetes'
$ws.Range('O9').Value = 'Provincia de Cautín'
$ws.Range('P9').Value = 667
$ws.Range('Q9').Value = 12
$ws.Range('R9').Value = 'Hortaliza'
